$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.061.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.69%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.017.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.86%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.607'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.80'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.60%  '
$ws.Range("E9").Value = '  -3.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0785'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("E11").Value = '  -5.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.316.26'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.95%  '
$ws.Range("E15").Value = '  -3.68%  '
$ws.Range("E16").Value = '  -3.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.016.81'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.011.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("E21").Value = '  -1.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = '  +1.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.44%  '
$ws.Range("E28").Value = '  -1.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.124'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.116'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.49'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.96%  '
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.24%  '
$ws.Range("E35").Value = '  -7.37%  '
$ws.Range("E36").Value = '  +1.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  -5.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.27'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0215'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.23%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.471.02'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '95.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0917'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.51%  '
$ws.Range("E44").Value = '  -4.26%  '
$ws.Range("E45").Value = '  -5.24%  '
$ws.Range("E46").Value = '  -5.90%  '
$ws.Range("E47").Value = '  -3.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.19'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.43%  '
$ws.Range("E49").Value = '  -1.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.203.54'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -13.71%  '
